# Add a new "Policies Audit" column to the database export header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell in column Y (25): "Policies Audit"
$ws.Cells.Item(1, 25).Value = "Policies Audit"

# Match the formatting of the other header cells (copy style from X1)
$ws.Cells.Item(1, 24).Copy()
$ws.Cells.Item(1, 25).PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Set widths for the new column (Y) and the following spare column (Z)
$ws.Columns.Item(25).ColumnWidth = 11.65
$ws.Columns.Item(26).ColumnWidth = 10.65

# Reset the active view back to the top-left of the sheet (A1)
$ws.Range("A1").Select()
